$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing value: G2 blue value changes from 0.1 to 0.3
$ws.Range("G2").Value = 0.3

# Row 6 is currently empty (just styled placeholder cells). Bring in the
# same formatting used by the rows above (string-style for A/B, numeric
# data-style for C:K) by copying formats from row 5, then populate values.
$ws.Range("A5:K5").Copy() | Out-Null
$ws.Range("A6:K6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A6").Value = "gruedisueli"
$ws.Range("B6").Value = "Gavin"
$ws.Range("C6").Value = 56
$ws.Range("D6").Value = 21
$ws.Range("E6").Value = 0.9
$ws.Range("F6").Value = 0.3
$ws.Range("G6").Value = 0.1
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = 0.3
$ws.Range("J6").Value = 0.35
$ws.Range("K6").Value = 0.39
